# Aggiornamento dati fino al 9 agosto 2021 (righe 329-343)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(329, 44403, 0, 3, 23.55712603062426),
    @(330, 44404, 0, 2, 15.70475068708284),
    @(331, 44405, 0, 2, 15.70475068708284),
    @(332, 44406, 1, 3, 23.55712603062426),
    @(333, 44407, 2, 4, 31.40950137416569),
    @(334, 44408, 2, 6, 47.11425206124853),
    @(335, 44409, 3, 8, 62.81900274833137),
    @(336, 44410, 1, 9, 70.67137809187278),
    @(337, 44411, 0, 9, 70.67137809187278),
    @(338, 44412, 1, 10, 78.52375343541422),
    @(339, 44413, 4, 13, 102.0808794660385),
    @(340, 44414, 1, 12, 94.22850412249706),
    @(341, 44415, 1, 11, 86.37612877895563),
    @(342, 44416, 5, 13, 102.0808794660385),
    @(343, 44417, 2, 14, 109.9332548095799)
)

foreach ($entry in $data) {
    $r = $entry[0]
    $ws.Cells.Item($r, 1).Value = $entry[1]
    $ws.Cells.Item($r, 2).Value = $entry[2]
    $ws.Cells.Item($r, 3).Value = $entry[3]
    $ws.Cells.Item($r, 4).Value = $entry[4]
}

# Match the date-style formatting (bold, thin border, centered, custom date
# number format) already used for column A by copying it from the row right
# above the newly-added block.
$ws.Cells.Item(328, 1).Copy()
$ws.Range("A329:A343").PasteSpecial(-4122)
